$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 2138
$ws.Range("C4").Value = 5560
$ws.Range("C5").Value = 6842
$ws.Range("C6").Value = 8981
$ws.Range("C7").Value = 11974
$ws.Range("C8").Value = 12402
$ws.Range("C9").Value = 14968
$ws.Range("C10").Value = 26087
$ws.Range("C11").Value = 40200
$ws.Range("C12").Value = 20528
$ws.Range("C13").Value = 20955
$ws.Range("C14").Value = 25231
$ws.Range("C15").Value = 24804
$ws.Range("C16").Value = 26515
$ws.Range("C17").Value = 30363
$ws.Range("C18").Value = 30791
$ws.Range("C19").Value = 32073
$ws.Range("C20").Value = 40199
$ws.Range("C21").Value = 41055
$ws.Range("C22").Value = 45759
$ws.Range("C23").Value = 42338
$ws.Range("C24").Value = 47042
$ws.Range("C25").Value = 60299
$ws.Range("C26").Value = 53457
$ws.Range("C27").Value = 51319
$ws.Range("C28").Value = 55167
$ws.Range("C29").Value = 53884
$ws.Range("C30").Value = 54312
$ws.Range("C31").Value = 56878
$ws.Range("C32").Value = 54740
$ws.Range("C33").Value = 63720
$ws.Range("C34").Value = 68425
$ws.Range("C35").Value = 85103
$ws.Range("C36").Value = 64575
$ws.Range("C37").Value = 75694
$ws.Range("C38").Value = 77833
$ws.Range("C39").Value = 78260
$ws.Range("C40").Value = 73556
$ws.Range("C41").Value = 82964
$ws.Range("C42").Value = 90235
$ws.Range("C43").Value = 118460
$ws.Range("C44").Value = 82110
$ws.Range("C45").Value = 90234
$ws.Range("C46").Value = 91090
$ws.Range("C47").Value = 98788
$ws.Range("C48").Value = 91945
$ws.Range("C49").Value = 97077
$ws.Range("C50").Value = 93229
$ws.Range("C51").Value = 91090
$ws.Range("C52").Value = 99215
$ws.Range("C53").Value = 98788
$ws.Range("C54").Value = 100498
$ws.Range("C55").Value = 101781
$ws.Range("C56").Value = 102209
$ws.Range("C57").Value = 107341
$ws.Range("C58").Value = 108197
$ws.Range("C59").Value = 108197
$ws.Range("C60").Value = 108624
$ws.Range("C61").Value = 116750
$ws.Range("C62").Value = 119743
$ws.Range("C63").Value = 111617
$ws.Range("C64").Value = 124875
$ws.Range("C65").Value = 122309
$ws.Range("C66").Value = 125302
$ws.Range("C67").Value = 124019
$ws.Range("C68").Value = 128296
$ws.Range("C69").Value = 134283
$ws.Range("C70").Value = 137704
$ws.Range("C71").Value = 131289
$ws.Range("C72").Value = 131717
$ws.Range("C73").Value = 139415
$ws.Range("C74").Value = 137704
$ws.Range("C75").Value = 143691
$ws.Range("C76").Value = 144118
$ws.Range("C77").Value = 144547
$ws.Range("C78").Value = 146685
$ws.Range("C79").Value = 145402
$ws.Range("C80").Value = 148823
$ws.Range("C81").Value = 202279
$ws.Range("C82").Value = 155666
$ws.Range("C83").Value = 148396
$ws.Range("C84").Value = 158660
$ws.Range("C85").Value = 159087
$ws.Range("C86").Value = 169778
$ws.Range("C87").Value = 176620
$ws.Range("C88").Value = 165074
$ws.Range("C89").Value = 237347
$ws.Range("C90").Value = 169778
$ws.Range("C91").Value = 163363
$ws.Range("C92").Value = 175766
$ws.Range("C93").Value = 169778
$ws.Range("C94").Value = 177476
$ws.Range("C95").Value = 174910
$ws.Range("C96").Value = 170633
$ws.Range("C97").Value = 181753
$ws.Range("C98").Value = 180042
$ws.Range("C99").Value = 217675
$ws.Range("C100").Value = 206556
$ws.Range("C101").Value = 193299
